$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update validation data value (E25): 4 -> 5
$ws.Range("E25").Value = 5

# Recalculate formulas so dependent totals/averages update
$excel.Calculate()

# Update the active selection to reflect the diff (G27)
$ws.Range("G27").Select()
